$wb = $excel.ActiveWorkbook

# ---------- Sheet 1: quality_comparison ----------
$ws1 = $wb.Worksheets.Item("quality_comparison")

# Give the top-right header cells (C1, D1) their box-closing borders.
# C1 gets top+bottom thin borders (matches border index 4 already defined in styles).
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item(8).Weight = 2   # xlEdgeTop
$c1.Borders.Item(9).Weight = 2   # xlEdgeBottom

# D1 gets top+bottom+right thin borders (matches border index 5 already defined in styles).
$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders.Item(8).Weight = 2    # xlEdgeTop
$d1.Borders.Item(10).Weight = 2   # xlEdgeRight
$d1.Borders.Item(9).Weight = 2    # xlEdgeBottom

# Anonymize "fedcore" column header label -> "approach"
$ws1.Range("C2").Value = "approach"

# Normalize "-0" change values to plain "0"
$ws1.Range("D4").Value = 0
$ws1.Range("D5").Value = 0

# ---------- Sheet 2: computational_comparison ----------
$ws2 = $wb.Worksheets.Item("computational_comparison")

$c1b = $ws2.Range("C1")
$c1b.ClearFormats()
$c1b.Borders.Item(8).Weight = 2
$c1b.Borders.Item(9).Weight = 2

$d1b = $ws2.Range("D1")
$d1b.ClearFormats()
$d1b.Borders.Item(8).Weight = 2
$d1b.Borders.Item(10).Weight = 2
$d1b.Borders.Item(9).Weight = 2

$f1b = $ws2.Range("F1")
$f1b.ClearFormats()
$f1b.Borders.Item(8).Weight = 2
$f1b.Borders.Item(9).Weight = 2

$g1b = $ws2.Range("G1")
$g1b.ClearFormats()
$g1b.Borders.Item(8).Weight = 2
$g1b.Borders.Item(10).Weight = 2
$g1b.Borders.Item(9).Weight = 2

# Anonymize "fedcore" column header labels -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5
$ws2.Range("G5").ClearContents()

Write-Host "edit applied"
